# Updated symbol list on Sat Dec 24 03:54:53 UTC 2022 with GitHub Actions
#
# This script reproduces the crypto price/ranking refresh captured in the
# target diff: most rows just get a refreshed Price (column D) value, while
# rows 10-18 are a coin that dropped out of the ranking ("FTXToken" was
# removed, everyone below it moved up one slot, and "One" was pulled back in
# at the bottom of that block) so B (Coin), C (Link), D (Price) and E
# (Volume(1h)) all shift for those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to push a literal (non-numeric-coerced) text value into a
# cell: building the value via a formula that returns a string, then pasting
# only the *value* back over the target, keeps it as plain text (matching the
# original inline-string price cells) instead of Excel auto-converting a
# numeric-looking string into a real number.
$scratch = "Z1"

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Text
    )
    $escaped = $Text.Replace("""", """""")
    $ws.Range($scratch).Formula = "=""" + $escaped + """"
    $ws.Range($scratch).Copy() | Out-Null
    $ws.Range($CellRef).PasteSpecial(-4163) | Out-Null
}

# --- Simple price-only refreshes -----------------------------------------
Set-TextValue "D2"  "245.89"
Set-TextValue "D3"  "22.15"
Set-TextValue "D4"  "5.358"
Set-TextValue "D5"  "0.05927"
Set-TextValue "D7"  "6.392"
Set-TextValue "D8"  "0.8128"
Set-TextValue "D9"  "0.9620"

# --- Rows 10-18: coin list shifted (One re-enters, rest shift down) ------
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D10" "0.01122"
$ws.Range("E10").Value = "9OneONEBestin24h"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1429"
$ws.Range("E11").Value = "10WazirXWRX"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07408"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03494"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D14" "0.03044"
$ws.Range("E14").Value = "13BitrueCoinBTR"

$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D15" "0.09399"
$ws.Range("E15").Value = "14BitMartTokenBMX"

$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D16" "3.999"
$ws.Range("E16").Value = "15MCDexMCB"

$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001598"
$ws.Range("E17").Value = "16BitForexTokenBF"

$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D18" "0.04812"
$ws.Range("E18").Value = "17CoinExTokenCET"

# --- More simple price-only refreshes -------------------------------------
Set-TextValue "D19" "0.006016"
Set-TextValue "D20" "0.004084"
Set-TextValue "D21" "0.0009874"
Set-TextValue "D22" "0.00009702"
Set-TextValue "D24" "2.164"
Set-TextValue "D40" "0.03939"
Set-TextValue "D41" "0.006472"
Set-TextValue "D42" "0.1072"
Set-TextValue "D43" "0.003001"
Set-TextValue "D44" "0.005479"

Set-TextValue "D47" "0.6601"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-TextValue "D48" "0.04639"
Set-TextValue "D49" "0.00002100"

# Clean up the scratch cell/clipboard state so nothing extra is left behind.
$ws.Range($scratch).Clear() | Out-Null
$excel.CutCopyMode = 0
